$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.147.35"
$ws.Range("E2").Value = "  -0.06%  "
$ws.Range("D3").Value = "2.052.27"
$ws.Range("E3").Value = "  -1.37%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "248.21"
$ws.Range("E5").Value = "  -2.58%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.664"
$ws.Range("E6").Value = "  -2.38%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.11"
$ws.Range("E7").Value = "  -6.92%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.382"
$ws.Range("E9").Value = "  -3.09%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0779"
$ws.Range("E10").Value = "  -2.88%  "
$ws.Range("E11").Value = "  -0.66%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.90"
$ws.Range("E12").Value = "  -2.12%  "
$ws.Range("D13").Value = "2.353.49"
$ws.Range("E13").Value = "  -1.24%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.834"
$ws.Range("E14").Value = "  +0.64%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.70"
$ws.Range("E15").Value = "  +2.19%  "
$ws.Range("D16").Value = "2.060.68"
$ws.Range("E16").Value = "  -1.01%  "
$ws.Range("E17").Value = "  +15.22%  "
$ws.Range("D18").Value = "37.210.50"
$ws.Range("E18").Value = "  +0.12%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "74.78"
$ws.Range("E19").Value = "  -0.12%  "
$ws.Range("D20").Value = "0.0₃0896"
$ws.Range("E20").Value = "  -3.83%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.34"
$ws.Range("E21").Value = "  -2.78%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "237.05"
$ws.Range("E22").Value = "  -1.80%  "
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("E24").Value = "  +1.54%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.17"
$ws.Range("E25").Value = "  -6.82%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "169.40"
$ws.Range("E26").Value = "  -0.56%  "
$ws.Range("E27").Value = "  +0.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.03"
$ws.Range("E28").Value = "  -2.17%  "
$ws.Range("E29").Value = "  -2.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.12"
$ws.Range("E30").Value = "  -0.44%  "
$ws.Range("E31").Value = "  -1.27%  "
$ws.Range("E32").Value = "  -4.00%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.48"
$ws.Range("E33").Value = "  -0.37%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0898"
$ws.Range("E34").Value = "  -2.00%  "
$ws.Range("E35").Value = "  +0.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.26"
$ws.Range("E36").Value = "  -2.28%  "
$ws.Range("E37").Value = "  -0.70%  "
$ws.Range("E38").Value = "  -3.13%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.23"
$ws.Range("E39").Value = "  +14.82%  "
$ws.Range("E40").Value = "  +15.01%  "
$ws.Range("E41").Value = "  -13.75%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0222"
$ws.Range("E42").Value = "  -3.09%  "
$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.18"
$ws.Range("E43").Value = "  -6.37%  "
$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.14"
$ws.Range("E44").Value = "  -3.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "95.77"
$ws.Range("E45").Value = "  -4.18%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.45"
$ws.Range("E46").Value = "  -2.41%  "
$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "1.275.39"
$ws.Range("E47").Value = "  -2.74%  "
$ws.Range("B48").Value = "MXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.86"
$ws.Range("E48").Value = "  -3.72%  "
$ws.Range("E49").Value = "  -2.16%  "
$ws.Range("D50").Value = "2.238.84"
$ws.Range("E50").Value = "  -1.26%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "43.74"
$ws.Range("E51").Value = "  -1.45%  "
